$wb = $excel.ActiveWorkbook

# Sheet ALC, row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 17990.334
$ws.Range("I86").Value = 1475
$ws.Range("K86").Value = 1475
$ws.Range("M86").Value = -352

# Sheet ALC, row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 17990.334
$ws.Range("I89").Value = 1475
$ws.Range("K89").Value = 7375
$ws.Range("M89").Value = -1759

# Sheet ALC, row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1142.5238
$ws.Range("J107").Value = 1841
$ws.Range("L107").Value = 1841
$ws.Range("N107").Value = -5681

# Sheet ALC, row 123
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -39800

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1913.0769
$ws.Range("I137").Value = 1686.7
$ws.Range("K137").Value = 5060.1
$ws.Range("M137").Value = -2510.1

# Sheet ARM, row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2814.8076
$ws.Range("I45").Value = 2506.3333
$ws.Range("J45").Value = 3235.4546
$ws.Range("K45").Value = 2506.3333
$ws.Range("L45").Value = 3235.4546
$ws.Range("M45").Value = -2129.3333
$ws.Range("N45").Value = -3989.4546

# Sheet ARM, row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1978.2188
$ws.Range("I61").Value = 1492.9642
$ws.Range("K61").Value = 1492.9642
$ws.Range("M61").Value = -1280.9642

# Sheet ARM, row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 100001260
$ws.Range("I97").Value = 1668.7142
$ws.Range("K97").Value = 1668.7142
$ws.Range("M97").Value = -1172.7142

# Sheet ARM, row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2575.5186
$ws.Range("I122").Value = 2195.7646
$ws.Range("K122").Value = 6587.293799999999
$ws.Range("M122").Value = -4137.293799999999

# Sheet ARM, row 123
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 35000
$ws.Range("J123").Value = 35000
$ws.Range("L123").Value = 35000
$ws.Range("N123").Value = -44800

# Sheet ARM, row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1978.2188
$ws.Range("I136").Value = 1492.9642
$ws.Range("K136").Value = 4478.892599999999
$ws.Range("M136").Value = -1928.892599999999

# Sheet BSM, row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2550
$ws.Range("I20").Value = 2600
$ws.Range("K20").Value = 2600
$ws.Range("M20").Value = -2353

# Sheet BSM, row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 851.3
$ws.Range("I94").Value = 666.1429000000001
$ws.Range("K94").Value = 666.1429000000001
$ws.Range("M94").Value = -215.1429000000001

# Sheet CRP, row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 758.5
$ws.Range("I16").Value = 705.8889
$ws.Range("J16").Value = 853.2
$ws.Range("K16").Value = 705.8889
$ws.Range("L16").Value = 853.2
$ws.Range("M16").Value = -418.8889
$ws.Range("N16").Value = -1427.2

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2903.7234
$ws.Range("I31").Value = 2241.5264
$ws.Range("K31").Value = 2241.5264
$ws.Range("M31").Value = -1946.5264

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2903.7234
$ws.Range("I34").Value = 2241.5264
$ws.Range("K34").Value = 2241.5264
$ws.Range("M34").Value = -2039.5264

# Sheet CRP, row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 18923.242
$ws.Range("I58").Value = 1615.0769
$ws.Range("J58").Value = 32986.125
$ws.Range("K58").Value = 1615.0769
$ws.Range("L58").Value = 32986.125
$ws.Range("M58").Value = -1412.0769
$ws.Range("N58").Value = -33392.125

# Sheet CRP, row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 758.5
$ws.Range("I113").Value = 705.8889
$ws.Range("J113").Value = 853.2
$ws.Range("K113").Value = 705.8889
$ws.Range("L113").Value = 853.2
$ws.Range("M113").Value = 1464.1111
$ws.Range("N113").Value = -5193.2

# Sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3144.0833
$ws.Range("I132").Value = 2475.2942
$ws.Range("K132").Value = 7425.882599999999
$ws.Range("M132").Value = -4895.882599999999

# Sheet CRP, row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1107.4375
$ws.Range("I134").Value = 702.1111
$ws.Range("K134").Value = 2106.3333
$ws.Range("M134").Value = 428.6667000000002

# Sheet CRP, row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 18923.242
$ws.Range("I136").Value = 1615.0769
$ws.Range("J136").Value = 32986.125
$ws.Range("K136").Value = 4845.2307
$ws.Range("L136").Value = 98958.375
$ws.Range("M136").Value = -2295.2307
$ws.Range("N136").Value = -104058.375

# Sheet CUL, row 46
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

# Sheet CUL, row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 550
$ws.Range("I122").Value = 317.16666
$ws.Range("J122").Value = 829.4
$ws.Range("K122").Value = 2854.49994
$ws.Range("L122").Value = 7464.599999999999
$ws.Range("M122").Value = -404.4999399999997
$ws.Range("N122").Value = -12364.6

# Sheet CUL, row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 740.71716
$ws.Range("J131").Value = 777.13336
$ws.Range("L131").Value = 2331.40008
$ws.Range("N131").Value = -12411.40008

# Sheet GSM, row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10436133
$ws.Range("J70").Value = 10436133
$ws.Range("L70").Value = 10436133
$ws.Range("N70").Value = -10436673

# Sheet GSM, row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 10436133
$ws.Range("J73").Value = 10436133
$ws.Range("L73").Value = 10436133
$ws.Range("N73").Value = -10438005

# Sheet GSM, row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4017.4
$ws.Range("I80").Value = 3410.3333
$ws.Range("J80").Value = 4358.875
$ws.Range("K80").Value = 3410.3333
$ws.Range("L80").Value = 4358.875
$ws.Range("M80").Value = -2412.3333
$ws.Range("N80").Value = -6354.875

# Sheet GSM, row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 4017.4
$ws.Range("I83").Value = 3410.3333
$ws.Range("J83").Value = 4358.875
$ws.Range("K83").Value = 17051.6665
$ws.Range("L83").Value = 21794.375
$ws.Range("M83").Value = -12059.6665
$ws.Range("N83").Value = -31778.375

# Sheet LTW, row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4199.875
$ws.Range("I22").Value = 3700
$ws.Range("J22").Value = 4499.8
$ws.Range("K22").Value = 3700
$ws.Range("L22").Value = 4499.8
$ws.Range("M22").Value = -3405
$ws.Range("N22").Value = -5089.8

# Sheet LTW, row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 4199.875
$ws.Range("I27").Value = 3700
$ws.Range("J27").Value = 4499.8
$ws.Range("K27").Value = 3700
$ws.Range("L27").Value = 4499.8
$ws.Range("M27").Value = -3593
$ws.Range("N27").Value = -4713.8

# Sheet LTW, row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 686.7646999999999
$ws.Range("I46").Value = 588.8889
$ws.Range("J46").Value = 796.875
$ws.Range("K46").Value = 588.8889
$ws.Range("L46").Value = 796.875
$ws.Range("M46").Value = -400.8889
$ws.Range("N46").Value = -1172.875

# Sheet LTW, row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3843.8333
$ws.Range("I61").Value = 1384.1538
$ws.Range("K61").Value = 1384.1538
$ws.Range("M61").Value = -1182.1538

# Sheet LTW, row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2380.5
$ws.Range("I68").Value = 2422.3333
$ws.Range("J68").Value = 2338.6667
$ws.Range("K68").Value = 2422.3333
$ws.Range("L68").Value = 2338.6667
$ws.Range("M68").Value = -1673.3333
$ws.Range("N68").Value = -3836.6667

# Sheet LTW, row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2380.5
$ws.Range("I71").Value = 2422.3333
$ws.Range("J71").Value = 2338.6667
$ws.Range("K71").Value = 12111.6665
$ws.Range("L71").Value = 11693.3335
$ws.Range("M71").Value = -8367.666499999999
$ws.Range("N71").Value = -19181.3335

# Sheet LTW, row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2263.2632
$ws.Range("I82").Value = 2535.6428
$ws.Range("J82").Value = 1500.6
$ws.Range("K82").Value = 2535.6428
$ws.Range("L82").Value = 1500.6
$ws.Range("M82").Value = -2174.6428
$ws.Range("N82").Value = -2222.6

# Sheet LTW, row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2263.2632
$ws.Range("I85").Value = 2535.6428
$ws.Range("J85").Value = 1500.6
$ws.Range("K85").Value = 2535.6428
$ws.Range("L85").Value = 1500.6
$ws.Range("M85").Value = -1287.6428
$ws.Range("N85").Value = -3996.6

# Sheet LTW, row 98
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H98").Value = 22000
$ws.Range("J98").Value = 22000
$ws.Range("L98").Value = 22000
$ws.Range("N98").Value = -27990

# Sheet LTW, row 108
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 19813
$ws.Range("J108").Value = 19813
$ws.Range("L108").Value = 19813
$ws.Range("N108").Value = -27493

# Sheet LTW, row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3843.8333
$ws.Range("I113").Value = 1384.1538
$ws.Range("K113").Value = 1384.1538
$ws.Range("M113").Value = 785.8462

# Sheet LTW, row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 936509.4
$ws.Range("I122").Value = 1309606.5
$ws.Range("K122").Value = 3928819.5
$ws.Range("M122").Value = -3926369.5

# Sheet LTW, row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1500
$ws.Range("I136").Value = 1500
$ws.Range("K136").Value = 4500
$ws.Range("M136").Value = -1950

# Sheet WVR, row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4166.5557
$ws.Range("I62").Value = 2500
$ws.Range("K62").Value = 2500
$ws.Range("M62").Value = -1876

# Sheet WVR, row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 4166.5557
$ws.Range("I65").Value = 2500
$ws.Range("K65").Value = 12500
$ws.Range("M65").Value = -9380

# Sheet WVR, row 101
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 16367.333
$ws.Range("J101").Value = 16367.333
$ws.Range("L101").Value = 16367.333
$ws.Range("N101").Value = -22857.333

# Sheet WVR, row 108
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 41000
$ws.Range("J108").Value = 41000
$ws.Range("L108").Value = 41000
$ws.Range("N108").Value = -48680

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1151.25
$ws.Range("I132").Value = 780.9375
$ws.Range("K132").Value = 2342.8125
$ws.Range("M132").Value = 187.1875
